$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '51.369.98'
$ws.Cells.Item(2, 5).Value = '  -0.55%  '
$ws.Cells.Item(3, 4).Value = '2.915.95'
$ws.Cells.Item(3, 5).Value = '  +0.74%  '
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Formula = "'363.49"
$ws.Cells.Item(5, 5).Value = '  +2.69%  '
$ws.Cells.Item(6, 4).Formula = "'104.60"
$ws.Cells.Item(6, 5).Value = '  -3.76%  '
$ws.Cells.Item(7, 5).Value = '  -2.97%  '
$ws.Cells.Item(8, 5).Value = '  -0.09%  '
$ws.Cells.Item(9, 5).Value = '  -5.00%  '
$ws.Cells.Item(10, 4).Formula = "'36.91"
$ws.Cells.Item(10, 5).Value = '  -4.52%  '
$ws.Cells.Item(11, 5).Value = '  +1.85%  '
$ws.Cells.Item(12, 4).Formula = "'0.0837"
$ws.Cells.Item(12, 5).Value = '  -3.32%  '
$ws.Cells.Item(13, 4).Formula = "'18.58"
$ws.Cells.Item(13, 5).Value = '  -4.04%  '
$ws.Cells.Item(14, 4).Value = '3.374.53'
$ws.Cells.Item(14, 5).Value = '  +0.47%  '
$ws.Cells.Item(15, 4).Formula = "'7.37"
$ws.Cells.Item(15, 5).Value = '  -3.93%  '
$ws.Cells.Item(16, 4).Value = '2.922.03'
$ws.Cells.Item(16, 5).Value = '  +0.63%  '
$ws.Cells.Item(17, 4).Formula = "'0.955"
$ws.Cells.Item(17, 5).Value = '  -1.48%  '
$ws.Cells.Item(18, 4).Value = '51.277.66'
$ws.Cells.Item(18, 5).Value = '  -0.69%  '
$ws.Cells.Item(19, 4).Formula = "'3.30"
$ws.Cells.Item(19, 5).Value = '  -1.90%  '
$ws.Cells.Item(20, 4).Formula = "'7.24"
$ws.Cells.Item(20, 5).Value = '  -3.23%  '
$ws.Cells.Item(21, 4).Formula = "'13.05"
$ws.Cells.Item(21, 5).Value = '  -4.97%  '
$ws.Cells.Item(22, 4).Value = '0.0₃0947'
$ws.Cells.Item(22, 5).Value = '  -2.50%  '
$ws.Cells.Item(23, 4).Formula = "'68.45"
$ws.Cells.Item(23, 5).Value = '  -2.38%  '
$ws.Cells.Item(24, 4).Formula = "'260.30"
$ws.Cells.Item(24, 5).Value = '  -2.48%  '
$ws.Cells.Item(25, 5).Value = '  -2.81%  '
$ws.Cells.Item(26, 5).Value = '  -4.51%  '
$ws.Cells.Item(27, 5).Value = '  +0.05%  '
$ws.Cells.Item(28, 4).Formula = "'26.14"
$ws.Cells.Item(28, 5).Value = '  -2.02%  '
$ws.Cells.Item(29, 4).Formula = "'7.27"
$ws.Cells.Item(29, 5).Value = '  -3.15%  '
$ws.Cells.Item(30, 5).Value = '  +4.17%  '
$ws.Cells.Item(31, 5).Value = '  -4.11%  '
$ws.Cells.Item(32, 4).Formula = "'6.17"
$ws.Cells.Item(32, 5).Value = '  +1.56%  '
$ws.Cells.Item(33, 5).Value = '  -2.08%  '
$ws.Cells.Item(34, 4).Formula = "'34.99"
$ws.Cells.Item(34, 5).Value = '  -5.95%  '
$ws.Cells.Item(35, 4).Formula = "'51.21"
$ws.Cells.Item(35, 5).Value = '  -1.50%  '
$ws.Cells.Item(36, 5).Value = '  +0.23%  '
$ws.Cells.Item(37, 4).Formula = "'0.0425"
$ws.Cells.Item(37, 5).Value = '  -2.95%  '
$ws.Cells.Item(38, 4).Formula = "'2.82"
$ws.Cells.Item(38, 5).Value = '  +5.19%  '
$ws.Cells.Item(39, 5).Value = '  -0.26%  '
$ws.Cells.Item(40, 5).Value = '  -5.76%  '
$ws.Cells.Item(41, 5).Value = '  -5.77%  '
$ws.Cells.Item(42, 5).Value = '  -3.75%  '
$ws.Cells.Item(43, 4).Formula = "'22.45"
$ws.Cells.Item(43, 5).Value = '  -0.48%  '
$ws.Cells.Item(44, 4).Formula = "'119.95"
$ws.Cells.Item(44, 5).Value = '  +1.13%  '
$ws.Cells.Item(45, 4).Formula = "'2.15"
$ws.Cells.Item(45, 5).Value = '  -0.99%  '
$ws.Cells.Item(46, 4).Value = '2.078.28'
$ws.Cells.Item(46, 5).Value = '  -1.81%  '
$ws.Cells.Item(47, 5).Value = '  -6.23%  '
$ws.Cells.Item(48, 4).Formula = "'2.27"
$ws.Cells.Item(48, 5).Value = '  -8.33%  '
$ws.Cells.Item(49, 4).Value = '3.211.21'
$ws.Cells.Item(49, 5).Value = '  +0.69%  '
$ws.Cells.Item(50, 5).Value = '  -4.05%  '
$ws.Cells.Item(51, 5).Value = '  -6.82%  '
